# Updating the models for NRG, PCSun and Ulmeni
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2-97) forward by 13 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 13
}

# Update the production values in column B for rows 28-49 with the new model output
$newB = @{
    28 = 6
    29 = 37
    30 = 116
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
}

foreach ($row in $newB.Keys) {
    $ws.Cells.Item($row, 2).Value2 = $newB[$row]
}
